$d = $word.ActiveDocument

# Color used for all highlighted metrics: RGB(0x2C, 0x3E, 0x50)
# RGB(r,g,b) = r + g*256 + b*65536
$metricColor = 44 + (62 * 256) + (80 * 65536)   # 0x2C3E50 -> 5258796

function Set-MetricBold($ParagraphIndex, $SearchText) {
    $range = $d.Paragraphs($ParagraphIndex).Range
    $found = $range.Find.Execute($SearchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        $range.Font.Bold = 1
        $range.Font.Color = $metricColor
    }
}

# "• Discovered systematic race coding errors ... accuracy from 23% to 64%"
Set-MetricBold 10 "23%"
Set-MetricBold 10 "64%"

# "• Utilized advanced sampling methods ... from ±4.2% to ±2.1%, increasing voter turnout prediction accuracy from 71% to 87%, ..."
Set-MetricBold 12 "±4.2%"
Set-MetricBold 12 "±2.1%"
Set-MetricBold 12 "71%"
Set-MetricBold 12 "87%"

# "• Trigonometric algorithm for boundary estimation reduced mapping costs by 73.5%, saving campaigns and organizations $4.7M and enabling smaller nonprofits to conduct analysis"
Set-MetricBold 13 "73.5%"
Set-MetricBold 13 "$4.7M"

# "• Built real-time FEC analysis systems ... valued over $2 trillion"
Set-MetricBold 14 "$2"

# "• Modernized legacy ETL processes by implementing dbt and PySpark workflows, reducing processing time by 57%"
Set-MetricBold 19 "57%"

# "• Algorithmic innovation: Pioneered trigonometric boundary estimation reducing mapping costs 73.5%"
Set-MetricBold 55 "73.5%"

# "• $4.7M savings enabled nonprofit access"
Set-MetricBold 56 "$4.7M"

# "• 178% accuracy improvement in racial classification algorithms"
Set-MetricBold 58 "178%"
